# Append the new trade row (row 3) to the ticker log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Y"

# B3 looks like a date ("09/02/20") but must stay plain text, like B2
# ("08/25/20") above it -- force a text format before assigning so Excel
# doesn't auto-convert it to a date serial, then drop back to the default
# (unstyled) cell style to match the rest of the data rows.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "09/02/20"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "BUY"
$ws.Range("D3").Value = 23
$ws.Range("E3").Value = 122
$ws.Range("F3").Value = -2806
$ws.Range("G3").Value = 122
$ws.Range("H3").Value = 2806
$ws.Range("I3").Value = 23

# J3 (REALIZED_PROFIT) is blank for this row -- still touch the cell so it
# exists in the sheet (matching the widened A1:J3 dimension) while leaving
# it visually/formula-wise empty and with the default style.
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = ""
$ws.Range("J3").Style = "Normal"
